# Auto-generated edit script to apply 'F' column (想去人数 / want-to-go count) updates
# across sheets 1-4 (展览, 演出, 本地生活, 全部类型) per the commit diff.
$wb = $excel.ActiveWorkbook

# Sheet 1 - 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 448  # F3: 447 -> 448
$ws.Cells.Item(5, 6).Value = 364  # F5: 362 -> 364
$ws.Cells.Item(6, 6).Value = 564  # F6: 560 -> 564
$ws.Cells.Item(9, 6).Value = 275  # F9: 274 -> 275
$ws.Cells.Item(10, 6).Value = 391  # F10: 388 -> 391
$ws.Cells.Item(12, 6).Value = 718  # F12: 712 -> 718
$ws.Cells.Item(13, 6).Value = 767  # F13: 765 -> 767
$ws.Cells.Item(14, 6).Value = 1  # F14: 0 -> 1
$ws.Cells.Item(15, 6).Value = 2  # F15: 0 -> 2
$ws.Cells.Item(16, 6).Value = 1521  # F16: 1520 -> 1521
$ws.Cells.Item(17, 6).Value = 1521  # F17: 1520 -> 1521
$ws.Cells.Item(20, 6).Value = 1357  # F20: 1356 -> 1357
$ws.Cells.Item(22, 6).Value = 341  # F22: 335 -> 341
$ws.Cells.Item(25, 6).Value = 106  # F25: 104 -> 106
$ws.Cells.Item(26, 6).Value = 6669  # F26: 6648 -> 6669
$ws.Cells.Item(27, 6).Value = 5036  # F27: 5019 -> 5036
$ws.Cells.Item(28, 6).Value = 5036  # F28: 5019 -> 5036
$ws.Cells.Item(32, 6).Value = 209  # F32: 207 -> 209
$ws.Cells.Item(35, 6).Value = 1297  # F35: 1294 -> 1297
$ws.Cells.Item(37, 6).Value = 253  # F37: 252 -> 253
$ws.Cells.Item(41, 6).Value = 255  # F41: 252 -> 255
$ws.Cells.Item(43, 6).Value = 150  # F43: 149 -> 150
$ws.Cells.Item(44, 6).Value = 63  # F44: 62 -> 63

# Sheet 2 - 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(10, 6).Value = 14  # F10: 13 -> 14
$ws.Cells.Item(18, 6).Value = 246  # F18: 242 -> 246

# Sheet 3 - 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 205  # F4: 203 -> 205
$ws.Cells.Item(5, 6).Value = 67  # F5: 62 -> 67

# Sheet 4 - 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 448  # F4: 447 -> 448
$ws.Cells.Item(7, 6).Value = 205  # F7: 203 -> 205
$ws.Cells.Item(8, 6).Value = 67  # F8: 62 -> 67
$ws.Cells.Item(9, 6).Value = 364  # F9: 362 -> 364
$ws.Cells.Item(10, 6).Value = 564  # F10: 560 -> 564
$ws.Cells.Item(14, 6).Value = 391  # F14: 388 -> 391
$ws.Cells.Item(16, 6).Value = 718  # F16: 712 -> 718
$ws.Cells.Item(17, 6).Value = 767  # F17: 765 -> 767
$ws.Cells.Item(18, 6).Value = 1521  # F18: 1520 -> 1521
$ws.Cells.Item(19, 6).Value = 1521  # F19: 1520 -> 1521
$ws.Cells.Item(22, 6).Value = 1357  # F22: 1356 -> 1357
$ws.Cells.Item(24, 6).Value = 341  # F24: 335 -> 341
$ws.Cells.Item(26, 6).Value = 106  # F26: 104 -> 106
$ws.Cells.Item(29, 6).Value = 6669  # F29: 6648 -> 6669
$ws.Cells.Item(30, 6).Value = 5036  # F30: 5019 -> 5036
$ws.Cells.Item(31, 6).Value = 5036  # F31: 5019 -> 5036
$ws.Cells.Item(33, 6).Value = 209  # F33: 207 -> 209
$ws.Cells.Item(34, 6).Value = 1297  # F34: 1294 -> 1297
$ws.Cells.Item(37, 6).Value = 253  # F37: 252 -> 253
$ws.Cells.Item(45, 6).Value = 255  # F45: 252 -> 255
$ws.Cells.Item(46, 6).Value = 150  # F46: 149 -> 150
$ws.Cells.Item(47, 6).Value = 63  # F47: 62 -> 63
$ws.Cells.Item(50, 6).Value = 246  # F50: 242 -> 246
